{"js": "// Expand the D-PAS project description: the intro sentence gains detail\n// about client systems (EPOS / Pulse) and the kinds of artifacts produced,\n// and the closing sentence about the tech stack switches from\n// \"multi-module project ... Flink, Kafka and Zookeeper in it. It is\n// developed using micro service architecture.\" to\n// \"multi-tenant application ... Kafka , Flink and Zookeeper in its\n// architecture.\" (Flink loses its bold emphasis in the process).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the D-PAS description paragraph by a distinctive text fragment\n// rather than a hard-coded index, so the script is resilient to unrelated\n// document changes.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text || \"\";\n  if (t.indexOf(\"is a server that is developed to process the insurance policies\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not locate the D-PAS description paragraph\");\n}\n\n// --- Edit 1: \"It takes requests form EPOS and processes to create\n// proposals, policies and claims. \" becomes a longer sentence describing\n// clients (EPOS/Pulse) and additional output types (reports etc.). ---\nlet hit1 = target.search(\n  \"It takes requests form EPOS and processes to create proposals, policies and claims. \",\n  { matchCase: true }\n);\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length === 0) {\n  throw new Error(\"Could not find the 'It takes requests form EPOS...' sentence\");\n}\nhit1.items[0].insertText(\n  \"It takes requests form clients like EPOS and Pulse and processes them to create proposals, policies, claims, reports etc. \",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Edit 2: \"...It is a multi-module project and uses \" becomes\n// \"...It is a multi-tenant application and uses Kafka , \" (Kafka moves in\n// front of Flink and stops being bold). ---\nlet hit2 = target.search(\n  \"amunda and work as per the workflow diagrams. It is a multi-module project and uses \",\n  { matchCase: true }\n);\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length === 0) {\n  throw new Error(\"Could not find the 'multi-module project and uses' sentence\");\n}\nhit2.items[0].insertText(\n  \"amunda and work as per the workflow diagrams. It is a multi-tenant application and uses Kafka , \",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Edit 3: \"Flink\" stops being bold (becomes regular/complex-script\n// weight only, i.e. bCs instead of b). ---\nlet hitFlink = target.search(\"Flink\", { matchCase: true });\nhitFlink.load(\"items\");\nawait context.sync();\nif (hitFlink.items.length === 0) {\n  throw new Error(\"Could not find 'Flink'\");\n}\nhitFlink.items[0].font.set({ bold: false, boldBidirectional: true });\nawait context.sync();\n\n// --- Edit 4: the rest of the sentence \", Kafka and Zookeeper in it. It is\n// developed using micro service architecture.\" is replaced by\n// \" and Zookeeper in its architecture. \" ---\nlet hit3 = target.search(\n  \", Kafka and Zookeeper in it. It is developed using micro service architecture.\",\n  { matchCase: true }\n);\nhit3.load(\"items\");\nawait context.sync();\nif (hit3.items.length === 0) {\n  throw new Error(\"Could not find the 'Kafka and Zookeeper in it...' tail sentence\");\n}\nhit3.items[0].insertText(\" and Zookeeper in its architecture. \", \"Replace\");\nawait context.sync();\n", "ps1": "# Expand the D-PAS project description: the intro sentence gains detail\n# about client systems (EPOS / Pulse) and the kinds of artifacts produced,\n# and the closing sentence about the tech stack switches from\n# \"multi-module project ... Flink, Kafka and Zookeeper in it. It is\n# developed using micro service architecture.\" to\n# \"multi-tenant application ... Kafka , Flink and Zookeeper in its\n# architecture.\" (Flink loses its bold emphasis in the process).\n\n$d = $word.ActiveDocument\n\n# Locate the D-PAS description paragraph by a distinctive text fragment\n# rather than a hard-coded index, so the script is resilient to unrelated\n# document changes.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*is a server that is developed to process the insurance policies*\") {\n        $target = $p\n        break\n    }\n}\nif ($null -eq $target) {\n    throw \"Could not locate the D-PAS description paragraph\"\n}\n$paraRange = $target.Range\n\n# --- Edit 1: \"It takes requests form EPOS and processes to create\n# proposals, policies and claims. \" becomes a longer sentence describing\n# clients (EPOS/Pulse) and additional output types (reports etc.). ---\n$r1 = $paraRange.Duplicate\n$found1 = $r1.Find.Execute(\"It takes requests form EPOS and processes to create proposals, policies and claims. \")\nif (-not $found1) {\n    throw \"Could not find the 'It takes requests form EPOS...' sentence\"\n}\n$r1.Text = \"It takes requests form clients like EPOS and Pulse and processes them to create proposals, policies, claims, reports etc. \"\n\n# --- Edit 2: \"...It is a multi-module project and uses \" becomes\n# \"...It is a multi-tenant application and uses Kafka , \" (Kafka moves in\n# front of Flink and stops being bold). ---\n$r2 = $target.Range.Duplicate\n$found2 = $r2.Find.Execute(\"amunda and work as per the workflow diagrams. It is a multi-module project and uses \")\nif (-not $found2) {\n    throw \"Could not find the 'multi-module project and uses' sentence\"\n}\n$r2.Text = \"amunda and work as per the workflow diagrams. It is a multi-tenant application and uses Kafka , \"\n\n# --- Edit 3: \"Flink\" stops being bold (becomes regular/complex-script\n# weight only, i.e. bCs instead of b). ---\n$r3 = $target.Range.Duplicate\n$found3 = $r3.Find.Execute(\"Flink\")\nif (-not $found3) {\n    throw \"Could not find 'Flink'\"\n}\n$r3.Font.Bold = 0\n$r3.Font.BoldBi = 1\n\n# --- Edit 4: the rest of the sentence \", Kafka and Zookeeper in it. It is\n# developed using micro service architecture.\" is replaced by\n# \" and Zookeeper in its architecture. \" ---\n$r4 = $target.Range.Duplicate\n$found4 = $r4.Find.Execute(\", Kafka and Zookeeper in it. It is developed using micro service architecture.\")\nif (-not $found4) {\n    throw \"Could not find the 'Kafka and Zookeeper in it...' tail sentence\"\n}\n$r4.Text = \" and Zookeeper in its architecture. \"\n"}
